$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update tags for LeetCode problem 57 "Insert Interval" (row 66, column C)
# "#array" -> "#array #核心"
$ws.Cells.Item(66, 3).Value = "#array #" + [char]0x6838 + [char]0x5fc3

# Add new row 67 for LeetCode problem 252 "Meeting Rooms"
$ws.Cells.Item(67, 1).Value = 252
$ws.Cells.Item(67, 2).Value = "Meeting Rooms"
$ws.Cells.Item(67, 3).Value = $ws.Cells.Item(64, 3).Value2
$ws.Cells.Item(67, 4).Value = "easy"
$ws.Cells.Item(67, 5).Value = 5
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 5
$ws.Cells.Item(67, 8).Value = 45847
$ws.Cells.Item(67, 9).Value = 45847

# Copy date formatting from the row above so H67/I67 reuse the existing date style
$ws.Range("H66:I66").Copy()
$ws.Range("H67:I67").PasteSpecial(-4122)

# Match row height used by similarly-sized rows (e.g. row 64)
$ws.Rows.Item(67).RowHeight = $ws.Rows.Item(64).RowHeight

# Update the active selection to reflect the new last row, like the source workbook
$ws.Range("H67:I67").Select()
